# Chocobo Profits workbook update (scheduled runner)
# Refreshes market-board derived columns (currentAveragePrice*, LevePrice*,
# LeveProfit*) for a handful of leves across the crafting-job sheets after
# new Universalis price data came in.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33: Glazed and Confused (Clear Glass Lens)
$ws.Range("H33").Value = 356.25
$ws.Range("I33").Value = 298.625
$ws.Range("J33").Value = 394.66666
$ws.Range("K33").Value = 298.625
$ws.Range("L33").Value = 394.66666
$ws.Range("M33").Value = -69.625
$ws.Range("N33").Value = -852.66666

# Row 101: Edge of the Arcane (Cunning Craftsman's Tea)
$ws.Range("H101").Value = 868.6667
$ws.Range("I101").Value = 262
$ws.Range("J101").Value = 1354
$ws.Range("K101").Value = 786
$ws.Range("L101").Value = 4062
$ws.Range("M101").Value = 836
$ws.Range("N101").Value = -7306

# Row 125: Body over Mind (Grade 5 Dexterity Alkahest)
$ws.Range("H125").Value = 2173.3333
$ws.Range("I125").Value = 1933.3334
$ws.Range("J125").Value = 2333.3333
$ws.Range("K125").Value = 17400.0006
$ws.Range("L125").Value = 20999.9997
$ws.Range("M125").Value = -14940.0006
$ws.Range("N125").Value = -25919.9997

# Row 127: Liquid Competence (Competent Craftsman's Draught)
$ws.Range("H127").Value = 2110.7778
$ws.Range("I127").Value = 1049.25
$ws.Range("J127").Value = 2960
$ws.Range("K127").Value = 3147.75
$ws.Range("L127").Value = 8880
$ws.Range("M127").Value = 1812.25
$ws.Range("N127").Value = -18800

# Row 132: Fast-forwarding Flora (Growth Formula Lambda)
$ws.Range("H132").Value = 114604.25
$ws.Range("I132").Value = 150250.11
$ws.Range("J132").Value = 7666.6665
$ws.Range("K132").Value = 450750.33
$ws.Range("L132").Value = 22999.9995
$ws.Range("M132").Value = -448220.33
$ws.Range("N132").Value = -28059.9995

# Row 135: For Tired Minds (Grade 1 Gemsap of Intelligence)
$ws.Range("H135").Value = 1507
$ws.Range("I135").Value = 1443.5
$ws.Range("K135").Value = 12991.5
$ws.Range("M135").Value = -10456.5

# Row 138: All-night Crafting (Cunning Craftsman's Tisane)
$ws.Range("H138").Value = 2412.7537
$ws.Range("I138").Value = 1346.2646
$ws.Range("J138").Value = 3448.7715
$ws.Range("K138").Value = 4038.7938
$ws.Range("L138").Value = 10346.3145
$ws.Range("M138").Value = 1101.2062
$ws.Range("N138").Value = -20626.3145


$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust (Steel Ingot)
$ws.Range("H32").Value = 4571.81
$ws.Range("I32").Value = 2602.972
$ws.Range("J32").Value = 9392.069
$ws.Range("K32").Value = 2602.972
$ws.Range("L32").Value = 9392.069
$ws.Range("M32").Value = -2315.972
$ws.Range("N32").Value = -9966.069

# Row 132: Don't Bore Me, Ore Me (Mountain Chromite Ingot)
$ws.Range("H132").Value = 2444.6538
$ws.Range("I132").Value = 1689.1
$ws.Range("J132").Value = 4963.1665
$ws.Range("K132").Value = 5067.299999999999
$ws.Range("L132").Value = 14889.4995
$ws.Range("M132").Value = -2537.299999999999
$ws.Range("N132").Value = -19949.4995


$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin (Adamantite Nugget)
$ws.Range("H86").Value = 1617
$ws.Range("I86").Value = 1587.5385
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 1587.5385
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -464.5385000000001
$ws.Range("N86").Value = -4246

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) (Adamantite Nugget)
$ws.Range("H89").Value = 1617
$ws.Range("I89").Value = 1587.5385
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 7937.692500000001
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -2321.692500000001
$ws.Range("N89").Value = -21232

# Row 94: High Steal (High Steel Nugget)
$ws.Range("H94").Value = 1223.871
$ws.Range("I94").Value = 1192.9231
$ws.Range("J94").Value = 1384.8
$ws.Range("K94").Value = 1192.9231
$ws.Range("L94").Value = 1384.8
$ws.Range("M94").Value = -741.9231
$ws.Range("N94").Value = -2286.8

# Row 97: File under Dull (High Steel File)
$ws.Range("H97").Value = 10172.75
$ws.Range("I97").Value = 1897
$ws.Range("J97").Value = 35000
$ws.Range("K97").Value = 1897
$ws.Range("L97").Value = 35000
$ws.Range("M97").Value = -906
$ws.Range("N97").Value = -36982

# Row 105: Ingot to Wing It (Molybdenum Ingot)
$ws.Range("H105").Value = 2607.95
$ws.Range("I105").Value = 2597.8975
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 2597.8975
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -850.8975
$ws.Range("N105").Value = -6494

# Row 134: Ruthenium Supremium (Ruthenium Ingot)
$ws.Range("H134").Value = 3024.152
$ws.Range("I134").Value = 1592.6945
$ws.Range("J134").Value = 8177.4
$ws.Range("K134").Value = 4778.083500000001
$ws.Range("L134").Value = 24532.2
$ws.Range("M134").Value = -2243.083500000001
$ws.Range("N134").Value = -29602.2


$ws = $wb.Worksheets.Item("CRP")
# Row 132: Hull Lotta Damage (Ginseng Lumber)
$ws.Range("H132").Value = 2670.0789
$ws.Range("I132").Value = 1311.85
$ws.Range("J132").Value = 4179.222
$ws.Range("K132").Value = 3935.55
$ws.Range("L132").Value = 12537.666
$ws.Range("M132").Value = -1405.55
$ws.Range("N132").Value = -17597.666


$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap (Maple Syrup)
$ws.Range("H5").Value = 1831.1818
$ws.Range("I5").Value = 499.17648
$ws.Range("J5").Value = 6360
$ws.Range("K5").Value = 1497.52944
$ws.Range("L5").Value = 19080
$ws.Range("M5").Value = -1385.52944
$ws.Range("N5").Value = -19304

# Row 113: Can't Eat Just One (Night Vinegar)
$ws.Range("H113").Value = 740.0244
$ws.Range("I113").Value = 644
$ws.Range("J113").Value = 972.0833
$ws.Range("K113").Value = 1932
$ws.Range("L113").Value = 2916.2499
$ws.Range("M113").Value = 238
$ws.Range("N113").Value = -7256.2499

# Row 135: Not-so-secret Ingredient (Royal Maple Syrup)
$ws.Range("H135").Value = 1831.1818
$ws.Range("I135").Value = 499.17648
$ws.Range("J135").Value = 6360
$ws.Range("K135").Value = 4492.58832
$ws.Range("L135").Value = 57240
$ws.Range("M135").Value = -1957.58832
$ws.Range("N135").Value = -62310


$ws = $wb.Worksheets.Item("GSM")
# Row 102: Put the Metal to the Peddle (Durium Ingot)
$ws.Range("H102").Value = 2438.2683
$ws.Range("I102").Value = 1894.3928
$ws.Range("J102").Value = 3609.6924
$ws.Range("K102").Value = 1894.3928
$ws.Range("L102").Value = 3609.6924
$ws.Range("M102").Value = -272.3928000000001
$ws.Range("N102").Value = -6853.6924

# Row 108: Satisfactory Sewing (Stonegold Needle)
$ws.Range("H108").Value = 53000
$ws.Range("J108").Value = 53000
$ws.Range("L108").Value = 53000
$ws.Range("N108").Value = -60680

# Row 122: Awarding Academic Excellence (Ametrine)
$ws.Range("H122").Value = 3158.25
$ws.Range("I122").Value = 2489.077
$ws.Range("J122").Value = 4401
$ws.Range("K122").Value = 7467.231000000001
$ws.Range("L122").Value = 13203
$ws.Range("M122").Value = -5017.231000000001
$ws.Range("N122").Value = -18103

# Row 132: On Board for Lar (Lar Ingot)
$ws.Range("H132").Value = 3400.9565
$ws.Range("I132").Value = 1757.25
$ws.Range("J132").Value = 5194.091
$ws.Range("K132").Value = 5271.75
$ws.Range("L132").Value = 15582.273
$ws.Range("M132").Value = -2741.75
$ws.Range("N132").Value = -20642.273


$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad (Toad Leather)
$ws.Range("H40").Value = 4523.3687
$ws.Range("I40").Value = 4060.1
$ws.Range("J40").Value = 6260.625
$ws.Range("K40").Value = 4060.1
$ws.Range("L40").Value = 6260.625
$ws.Range("M40").Value = -3924.1
$ws.Range("N40").Value = -6532.625

# Row 122: Hell on Leather (Gaja Leather)
$ws.Range("H122").Value = 3269.6667
$ws.Range("I122").Value = 2900.6365
$ws.Range("J122").Value = 4622.778
$ws.Range("K122").Value = 8701.9095
$ws.Range("L122").Value = 13868.334
$ws.Range("M122").Value = -6251.9095
$ws.Range("N122").Value = -18768.334

# Row 132: Tenets of Tanning (Silver Lobo Leather)
$ws.Range("H132").Value = 7829.364
$ws.Range("I132").Value = 2867.6667
$ws.Range("J132").Value = 8612.789000000001
$ws.Range("K132").Value = 8603.000100000001
$ws.Range("L132").Value = 25838.367
$ws.Range("M132").Value = -6073.000100000001
$ws.Range("N132").Value = -30898.367


$ws = $wb.Worksheets.Item("WVR")
# Row 122: Heavy Armoire (Dark Hempen Cloth)
$ws.Range("H122").Value = 2592.7317
$ws.Range("I122").Value = 1922.7307
$ws.Range("J122").Value = 3754.0667
$ws.Range("K122").Value = 5768.1921
$ws.Range("L122").Value = 11262.2001
$ws.Range("M122").Value = -3318.1921
$ws.Range("N122").Value = -16162.2001

# Row 132: Comfy Cabins (Snow Cotton Cloth)
$ws.Range("H132").Value = 7753702.5
$ws.Range("I132").Value = 754.6087
$ws.Range("J132").Value = 16669592
$ws.Range("K132").Value = 2263.8261
$ws.Range("L132").Value = 50008776
$ws.Range("M132").Value = 266.1738999999998
$ws.Range("N132").Value = -50013836

# Row 136: Weaving the Envelope (Sarcenet Cloth)
$ws.Range("H136").Value = 3606.5715
$ws.Range("I136").Value = 758.3
$ws.Range("K136").Value = 2274.9
$ws.Range("M136").Value = 275.1000000000004
